$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2375
$ws.Range("J40").Value = 2428.5715
$ws.Range("L40").Value = 2428.5715
$ws.Range("N40").Value = -2778.5715
$ws.Range("H48").Value = 792.8570999999999
$ws.Range("I48").Value = 500
$ws.Range("J48").Value = 841.6667
$ws.Range("K48").Value = 1500
$ws.Range("L48").Value = 2525.0001
$ws.Range("M48").Value = -1208
$ws.Range("N48").Value = -3109.0001
$ws.Range("H56").Value = 792.8570999999999
$ws.Range("I56").Value = 500
$ws.Range("J56").Value = 841.6667
$ws.Range("K56").Value = 1500
$ws.Range("L56").Value = 2525.0001
$ws.Range("M56").Value = -966
$ws.Range("N56").Value = -3593.0001
$ws.Range("H97").Value = 2166.6667
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 9000
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -9992
$ws.Range("H101").Value = 1989.3334
$ws.Range("I101").Value = 484
$ws.Range("J101").Value = 5000
$ws.Range("K101").Value = 1452
$ws.Range("L101").Value = 15000
$ws.Range("M101").Value = 170
$ws.Range("N101").Value = -18244
$ws.Range("H116").Value = 4000
$ws.Range("H137").Value = 4284.4546
$ws.Range("I137").Value = 4132.5
$ws.Range("J137").Value = 4550.375
$ws.Range("K137").Value = 12397.5
$ws.Range("L137").Value = 13651.125
$ws.Range("M137").Value = -9847.5
$ws.Range("N137").Value = -18751.125
$ws.Range("H138").Value = 2236.543
$ws.Range("I138").Value = 1901.96
$ws.Range("J138").Value = 3073
$ws.Range("K138").Value = 5705.88
$ws.Range("L138").Value = 9219
$ws.Range("M138").Value = -565.8800000000001
$ws.Range("N138").Value = -19499

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1248.6
$ws.Range("I97").Value = 1212.5
$ws.Range("K97").Value = 1212.5
$ws.Range("M97").Value = -716.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 505
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 505
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 505
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -731
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 3078.2
$ws.Range("I134").Value = 2813.6428
$ws.Range("J134").Value = 3414.9092
$ws.Range("K134").Value = 8440.928400000001
$ws.Range("L134").Value = 10244.7276
$ws.Range("M134").Value = -5905.928400000001
$ws.Range("N134").Value = -15314.7276

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5051.35
$ws.Range("I31").Value = 1241.9706
$ws.Range("J31").Value = 10032.846
$ws.Range("K31").Value = 1241.9706
$ws.Range("L31").Value = 10032.846
$ws.Range("M31").Value = -946.9706000000001
$ws.Range("N31").Value = -10622.846
$ws.Range("H34").Value = 5051.35
$ws.Range("I34").Value = 1241.9706
$ws.Range("J34").Value = 10032.846
$ws.Range("K34").Value = 1241.9706
$ws.Range("L34").Value = 10032.846
$ws.Range("M34").Value = -1039.9706
$ws.Range("N34").Value = -10436.846
$ws.Range("H74").Value = 23635.092
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 25498.6
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 25498.6
$ws.Range("M74").Value = -4126
$ws.Range("N74").Value = -27246.6
$ws.Range("H77").Value = 23635.092
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 25498.6
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 76495.79999999999
$ws.Range("M77").Value = -10632
$ws.Range("N77").Value = -85231.79999999999
$ws.Range("H92").Value = 40719.8
$ws.Range("J92").Value = 40719.8
$ws.Range("L92").Value = 40719.8
$ws.Range("N92").Value = -45711.8
$ws.Range("H132").Value = 3402776.5
$ws.Range("I132").Value = 1094.3948
$ws.Range("K132").Value = 3283.1844
$ws.Range("M132").Value = -753.1844000000001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5507.5
$ws.Range("I56").Value = 5507.5
$ws.Range("K56").Value = 5507.5
$ws.Range("M56").Value = -4977.5
$ws.Range("H121").Value = 1142.4419
$ws.Range("I121").Value = 209.7
$ws.Range("J121").Value = 1425.091
$ws.Range("K121").Value = 629.0999999999999
$ws.Range("L121").Value = 4275.272999999999
$ws.Range("M121").Value = 680.9000000000001
$ws.Range("N121").Value = -6895.272999999999
$ws.Range("H129").Value = 1706.5264
$ws.Range("J129").Value = 2440.3333
$ws.Range("L129").Value = 7320.999899999999
$ws.Range("N129").Value = -17320.9999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1480.2222
$ws.Range("I31").Value = 1480.2222
$ws.Range("K31").Value = 1480.2222
$ws.Range("M31").Value = -1188.2222
$ws.Range("H37").Value = 1480.2222
$ws.Range("I37").Value = 1480.2222
$ws.Range("K37").Value = 1480.2222
$ws.Range("M37").Value = -1203.2222
$ws.Range("H46").Value = 4141.952
$ws.Range("I46").Value = 1994
$ws.Range("J46").Value = 4499.9443
$ws.Range("K46").Value = 1994
$ws.Range("L46").Value = 4499.9443
$ws.Range("M46").Value = -1838
$ws.Range("N46").Value = -4811.9443
$ws.Range("H103").Value = 55000
$ws.Range("J103").Value = 55000
$ws.Range("L103").Value = 55000
$ws.Range("N103").Value = -57344

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 71431300
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 71431300
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H122").Value = 2270.889
$ws.Range("I122").Value = 846.6667
$ws.Range("J122").Value = 2983
$ws.Range("K122").Value = 2540.0001
$ws.Range("L122").Value = 8949
$ws.Range("M122").Value = -90.0001000000002
$ws.Range("N122").Value = -13849

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 91666.664
$ws.Range("J76").Value = 91666.664
$ws.Range("L76").Value = 91666.664
$ws.Range("N76").Value = -92296.664
$ws.Range("H79").Value = 91666.664
$ws.Range("J79").Value = 91666.664
$ws.Range("L79").Value = 91666.664
$ws.Range("N79").Value = -93850.664
$ws.Range("H100").Value = 1339.1666
$ws.Range("I100").Value = 1333.8572
$ws.Range("J100").Value = 1357.75
$ws.Range("K100").Value = 2667.7144
$ws.Range("L100").Value = 2715.5
$ws.Range("M100").Value = -2126.7144
$ws.Range("N100").Value = -3797.5
$ws.Range("H136").Value = 1827.5
$ws.Range("I136").Value = 1337.5476
$ws.Range("J136").Value = 3542.3333
$ws.Range("K136").Value = 4012.642800000001
$ws.Range("L136").Value = 10626.9999
$ws.Range("M136").Value = -1462.642800000001
$ws.Range("N136").Value = -15726.9999
